$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Feed Retrofit Overview" sheet - HCRO Time Standard retrofit tracking grid
# Mark several feeds' checklist columns as "In Progress" / "x" (done)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Feed Retrofit Overview")

function Set-InProgress($addr) {
    $ws.Range("D4").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = "In Progress"
}

function Set-Done($addr) {
    $ws.Range("C7").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = "x"
}

# Row 4 - 5C4-003
Set-InProgress "C4"
Set-Done "I4"
Set-InProgress "J4"

# Row 5 - 5C4-004
Set-InProgress "C5"
Set-Done "I5"
Set-InProgress "J5"

# Row 9 - 5C4-008
Set-InProgress "C9"
Set-Done "E9"
Set-Done "H9"
Set-Done "I9"
Set-InProgress "J9"

# Row 11 - 5C4-010
Set-InProgress "C11"
Set-Done "I11"

# Row 12 - 5C4-011
Set-InProgress "C12"
Set-Done "I12"
Set-InProgress "J12"

# Row 15 - 5C4-014
Set-InProgress "C15"
Set-Done "E15"
Set-Done "I15"
Set-InProgress "J15"

# Print scaled to 69% on this sheet
$ws.PageSetup.Zoom = $false
$ws.PageSetup.Scale = 69

# Last cell the author had selected before saving
$ws.Range("E29").Select()

# ---------------------------------------------------------------------------
# Sheets that never had an explicit page setup yet: materialize Excel's
# default (A4, landscape) the same way opening Page Setup / Print Preview
# does.
# ---------------------------------------------------------------------------
foreach ($name in @("Parts at SRI", "Parts at Minex", "Parts at SSL", "PAX Boxes")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.PageSetup.PaperSize = 9
    $sheet.PageSetup.Orientation = 2
}

# ---------------------------------------------------------------------------
# Restore the window position recorded the last time the workbook was saved
# ---------------------------------------------------------------------------
$excel.Left = 41060
$excel.Top = 2900
